# Weekly update: insert 3 new rows of data (newest week, date 44511)
# at the top of the data block (row 369), pushing the existing rows
# 369:465 down to 372:468.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 369:465 down by 3 rows (-> 372:468)
$ws.Rows("369:371").Insert()

# Fixed (constant) column values shared by every row in this dataset
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$tipo      = "Fruta"
$productoId = 100101
$producto  = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad  = "Sin especificar"
$unidad    = "`$/bandeja 7 kilos"
$origen    = "Provincia de Melipilla"
$kgUnidad  = 7

# New rows (newest price week) to populate at 369, 370, 371
# columns: D, L, M, N, O, P, S
$newRows = @(
    @{ Row = 369; Fecha = 44511; Calidad = "Especial"; Volumen = 400; PrecioMin = 12000; PrecioMax = 12500; PrecioProm = 12250; PrecioKg = 1750 },
    @{ Row = 370; Fecha = 44511; Calidad = "Primera";  Volumen = 400; PrecioMin = 10000; PrecioMax = 10500; PrecioProm = 10250; PrecioKg = 1464 },
    @{ Row = 371; Fecha = 44511; Calidad = "Segunda";  Volumen = 360; PrecioMin = 8000;  PrecioMax = 8500;  PrecioProm = 8250;  PrecioKg = 1179 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $rowData.Fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $rowData.Calidad
    $ws.Cells.Item($r, 13).Value = $rowData.Volumen
    $ws.Cells.Item($r, 14).Value = $rowData.PrecioMin
    $ws.Cells.Item($r, 15).Value = $rowData.PrecioMax
    $ws.Cells.Item($r, 16).Value = $rowData.PrecioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $rowData.PrecioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
